$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Text number format to the whole Student ID column so the
# leading zeros in the ID values are preserved.
$ws.Range("A1:A3").NumberFormat = "@"

# New student row.
$ws.Range("A3").Value = "0012800586"
$ws.Range("B3").Value = "Juanna"

# Column widths (character units -> stored width includes ~0.8333 padding).
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 20.166666666666668

# Move/activate the selection like in the saved workbook.
$ws.Range("A10").Select()

$wb.Save()
